# Todo list workbook: persist newly added todo items.
# The app now uses local storage (client-side) for persistence and the
# exported worksheet reflects the up-to-date todo list, which now
# includes three additional "this is a task" items appended below the
# existing rows (drag-and-drop / local storage persistence fix).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing rows 1-4 (header + 3 original todos) are unchanged.
# Append the new todo rows at 5, 6, 7.

$ws.Range("A5").Value = "c3a18998-4432-473b-950a-7943003c0a6e"
$ws.Range("B5").Value = "this is a task"
$ws.Range("C5").Value = $false

$ws.Range("A6").Value = "80fc3588-1dcf-4e4e-8d3b-4520758457f9"
$ws.Range("B6").Value = "this is a task"
$ws.Range("C6").Value = $false

$ws.Range("A7").Value = "5d8abbcb-a305-487a-9872-8f8f75cd31e4"
$ws.Range("B7").Value = "this is a task"
$ws.Range("C7").Value = $false
